# First working edition with WinPatcher.
# Reorganise the benchmark sheet to make room for x86 results alongside
# the existing x64 results, and add the (currently empty) x86 section
# headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Step 1: move the existing "TortureTest x64" block (rows 10-14)
#     down to rows 30-34 to make room for the new x86 sections.
#     Capture values first, then clear the old locations, then write
#     the new locations - this avoids row 11 (shared data row) being
#     overwritten before we read it.

$tortureTitle = $ws.Range("A10").Value2
$row11 = @($ws.Range("B11").Value2, $ws.Range("C11").Value2, $ws.Range("D11").Value2, $ws.Range("E11").Value2, $ws.Range("F11").Value2)
$a12 = $ws.Range("A12").Value2
$f12 = $ws.Range("F12").Value2
$a13 = $ws.Range("A13").Value2
$f13 = $ws.Range("F13").Value2
$a14 = $ws.Range("A14").Value2
$f14 = $ws.Range("F14").Value2

# --- Step 2: move the existing "SpeedTest x64" block (rows 1-5) down
#     to rows 10-14 (same shape as the TortureTest block above).

$speedTitle = $ws.Range("A1").Value2
$row2 = @($ws.Range("B2").Value2, $ws.Range("C2").Value2, $ws.Range("D2").Value2, $ws.Range("E2").Value2, $ws.Range("F2").Value2)
$a3 = $ws.Range("A3").Value2
$f3 = $ws.Range("F3").Value2
$a4 = $ws.Range("A4").Value2
$f4 = $ws.Range("F4").Value2
$a5 = $ws.Range("A5").Value2
$f5 = $ws.Range("F5").Value2

# Clear out the old ranges before re-writing the sheet so nothing is
# left behind at the old row numbers.
$ws.Range("A1:G14").ClearContents()

# --- Write the SpeedTest x64 block at its new home: rows 10-14.
$ws.Range("A10").Value2 = $speedTitle
$ws.Range("B11").Value2 = $row2[0]
$ws.Range("C11").Value2 = $row2[1]
$ws.Range("D11").Value2 = $row2[2]
$ws.Range("E11").Value2 = $row2[3]
$ws.Range("F11").Value2 = $row2[4]
$ws.Range("A12").Value2 = $a3
$ws.Range("F12").Value2 = $f3
$ws.Range("A13").Value2 = $a4
$ws.Range("F13").Value2 = $f4
$ws.Range("G13").Formula = "=F13/F`$12"
$ws.Range("A14").Value2 = $a5
$ws.Range("F14").Value2 = $f5
$ws.Range("G14").Formula = "=F14/F`$12"

# --- Write the TortureTest x64 block at its new home: rows 30-34.
$ws.Range("A30").Value2 = $tortureTitle
$ws.Range("B31").Value2 = $row11[0]
$ws.Range("C31").Value2 = $row11[1]
$ws.Range("D31").Value2 = $row11[2]
$ws.Range("E31").Value2 = $row11[3]
$ws.Range("F31").Value2 = $row11[4]
$ws.Range("A32").Value2 = $a12
$ws.Range("F32").Value2 = $f12
$ws.Range("A33").Value2 = $a13
$ws.Range("F33").Value2 = $f13
$ws.Range("G33").Formula = "=F33/F`$32"
$ws.Range("A34").Value2 = $a14
$ws.Range("F34").Value2 = $f14
$ws.Range("G34").Formula = "=F34/F`$32"

# --- New section headers for the not-yet-populated x86 runs.
$ws.Range("A1").Value2 = "SpeedTest Windows XP on x86 quad-core (Intel 2.66Ghz Q6700)"
$ws.Range("A20").Value2 = "TortureTest Windows XP on x86 quad-core (Intel 2.66Ghz Q6700)"

# --- Restore selection to match the author's saved view.
$ws.Range("B7").Select()
